# Update EC database: reorder RUBEN BELTRAN ZABALA rows to the top (ascending
# period order), move ELENA CLAUDIA ROBLES CUESTA's rows after them, and drop
# the WILBER DIAZ ALTAMIRANDA row entirely (commit: "Actualiza base de datos
# EC y agrega parte 1 de nuevos estado de cuenta").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite the 7 data rows that survive (16-22) with their new contents ---

# Row 16: RUBEN BELTRAN ZABALA, periodo 1601
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "73093935"
$ws.Range("D16").Value = "RUBEN BELTRAN ZABALA"
$ws.Range("E16").Value = "1601"
$ws.Range("F16").Value = 25800
$ws.Range("G16").Value = 877803

# Row 17: RUBEN BELTRAN ZABALA, periodo 1602
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73093935"
$ws.Range("D17").Value = "RUBEN BELTRAN ZABALA"
$ws.Range("E17").Value = "1602"
$ws.Range("F17").Value = 25800
$ws.Range("G17").Value = 877803

# Row 18: RUBEN BELTRAN ZABALA, periodo 1603
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "73093935"
$ws.Range("D18").Value = "RUBEN BELTRAN ZABALA"
$ws.Range("E18").Value = "1603"
$ws.Range("F18").Value = 25800
$ws.Range("G18").Value = 877803

# Row 19: RUBEN BELTRAN ZABALA, periodo 1604
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73093935"
$ws.Range("D19").Value = "RUBEN BELTRAN ZABALA"
$ws.Range("E19").Value = "1604"
$ws.Range("F19").Value = 25800
$ws.Range("G19").Value = 877803

# Row 20: RUBEN BELTRAN ZABALA, periodo 1605
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "73093935"
$ws.Range("D20").Value = "RUBEN BELTRAN ZABALA"
$ws.Range("E20").Value = "1605"
$ws.Range("F20").Value = 25800
$ws.Range("G20").Value = 877803

# Row 21: ELENA CLAUDIA ROBLES CUESTA, periodo 2003
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "45486547"
$ws.Range("D21").Value = "ELENA CLAUDIA ROBLES CUESTA"
$ws.Range("E21").Value = "2003"
$ws.Range("F21").Value = 48000
$ws.Range("G21").Value = 1200000

# Row 22: RUBEN BELTRAN ZABALA, periodo 2003
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "73093935"
$ws.Range("D22").Value = "RUBEN BELTRAN ZABALA"
$ws.Range("E22").Value = "2003"
$ws.Range("F22").Value = 35112
$ws.Range("G22").Value = 877803

# --- Drop row 23 (WILBER DIAZ ALTAMIRANDA) entirely; this shifts the
#     signature footer rows (28/29 -> 27/28) up automatically ---
$ws.Rows("23").Delete()

# --- Update the summary figures above the table ---
$ws.Range("E11").Value = 212112   # VALOR MORA total
$ws.Range("C13").Value = 2        # Cant. Trabajadores
